$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title: "Support the Nagoya Protocol" -> "support the Nagoya Protocol"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Core Principles for an Online Permit and Monitoring System to Support the Nagoya Protocol", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Core Principles for an Online Permit and Monitoring System to support the Nagoya Protocol", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "...relevant to the terms and conditions of the permit." -> "...permit and associated MAT."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("relevant to the terms and conditions of the permit. Making it easier", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "relevant to the terms and conditions of the permit and associated MAT. Making it easier", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from the "online permit system can facilitate"
#    paragraph to immediately after "Access and Benefit Sharing Related Principles".
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$headingRange = $d.Content
$headingRange.Find.Execute("Access and Benefit Sharing Related Principles") | Out-Null
$endPos = $headingRange.End

# Degenerate (zero-length) ranges placed exactly at a paragraph-end boundary
# confuse Bookmarks.Add in this host, so insert a throwaway character, bookmark
# around it (a real, non-empty range), then delete the character's text via
# the bookmark's own Range -- this leaves a correctly positioned zero-width
# bookmark with no surrounding side effects.
$tmp = $d.Range($endPos, $endPos)
$tmp.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $tmp) | Out-Null
$bmRange = $d.Bookmarks("_GoBack").Range
$bmRange.Text = ""

# ---------------------------------------------------------------------------
# 4) Add a tab-stop override to the first "A Single System" list paragraph
#    (numId 3, ilvl 0): clear the 0pt tab and add one at 24pt (480 twips).
# ---------------------------------------------------------------------------
$singleSystem = $d.Content
$singleSystem.Find.Execute("A Single System") | Out-Null
$singleSystemPara = $singleSystem.Paragraphs(1)
$singleSystemTabs = $singleSystemPara.Range.ParagraphFormat.TabStops
$singleSystemTabs.ClearAll()
$singleSystemTabs.Add(24, 0) | Out-Null

# ---------------------------------------------------------------------------
# 5) Remove the tab-stop override from the six "Compact" list paragraphs
#    (numId 11) under "The online permit system can facilitate information...".
# ---------------------------------------------------------------------------
$compactAnchors = @(
    "Numbers of permits granted by type.",
    "Organisations/Companies involved.",
    "Funding bodies involved.",
    "Countries involved.",
    "Publications/patent applications or products arising.",
    "Provision of information to the ABS Clearing House Mechanism."
)

foreach ($anchor in $compactAnchors) {
    $r = $d.Content
    $r.Find.Execute($anchor) | Out-Null
    $p = $r.Paragraphs(1)
    $p.Range.ParagraphFormat.TabStops.ClearAll()
}

Write-Output "done"
